# K562_A549_PeakN.xlsx update
# - Sheet1: add a new replicate row (H3K4me3_K562_E4) as row 7, shifting the
#   rest of the table down (Excel auto-adjusts the D-column ratio formulas).
# - Add a new "Sheet2" that consolidates the Sample/PeakN.after.tidy table
#   (minus the repeated header row) under new Exp_Name / Peak_Number headers,
#   and make it the active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: insert the new H3K4me3_K562_E4 replicate as row 7 -------------
$ws1.Rows("7:7").Insert()
$ws1.Range("A7").Value = "H3K4me3_K562_E4"
$ws1.Range("B7").Value = 43357
# The inserted row copied the style of the row above into column D; the
# source table has no ratio formula for this row, so drop that stray cell.
$ws1.Range("D7").Clear()

# Update the visible selection/scroll position left behind by the edit.
$ws1.Range("A17:B27").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1

# --- Add Sheet2 (placed after Sheet1) and populate the consolidated table --
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Exp_Name"
$ws2.Range("B1").Value = "Peak_Number"

$names = @(
  "H3K4me1_K562_E1", "H3K4me1_K562_E2", "H3K4me3_K562_E1", "H3K4me3_K562_E2",
  "H3K4me3_K562_E3", "H3K4me3_K562_E4", "H3K9ac_K562_E1", "H3K9ac_K562_E2",
  "H3K27ac_K562_E1", "H3K27me3_K562_E1", "H3K27me3_K562_E2", "H3K36me3_K562_E1",
  "H3K36me3_K562_E2", "H3K4me1_A549_E1", "H3K4me1_A549_E2", "H3K4me3_A549_E1",
  "H3K4me3_A549_E2", "H3K9ac_A549_E1", "H3K27ac_A549_E1", "H3K27ac_A549_E2",
  "H3K27me3_A549_E1", "H3K27me3_A549_E2", "H3K36me3_A549_E1", "H3K36me3_A549_E2"
)
$values = @(
  112834, 95383, 21146, 23132,
  29544, 43357, 120499, 38252,
  51343, 145139, 119855, 74176,
  161791, 124322, 130082, 27884,
  32532, 47957, 75735, 47542,
  139963, 172278, 136039, 91433
)

for ($i = 0; $i -lt $names.Count; $i++) {
  $row = $i + 2
  $ws2.Cells.Item($row, 1).Value = $names[$i]
  $ws2.Cells.Item($row, 2).Value = $values[$i]
}

# Sheet2 becomes the active/visible tab (matches the saved workbook state).
$ws2.Activate()
$ws2.Range("G17").Select()
